$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2241.56
$ws.Range("J19").Value = 1644.8462
$ws.Range("L19").Value = 1644.8462
$ws.Range("N19").Value = -1994.8462
$ws.Range("H137").Value = 6998.75
$ws.Range("I137").Value = 10730.538
$ws.Range("J137").Value = 3764.5334
$ws.Range("K137").Value = 32191.614
$ws.Range("L137").Value = 11293.6002
$ws.Range("M137").Value = -29641.614
$ws.Range("N137").Value = -16393.6002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5340.7256
$ws.Range("I32").Value = 5387.295
$ws.Range("J32").Value = 2500
$ws.Range("K32").Value = 5387.295
$ws.Range("L32").Value = 2500
$ws.Range("M32").Value = -5100.295
$ws.Range("N32").Value = -3074
$ws.Range("H45").Value = 7212.826
$ws.Range("I45").Value = 8175.647
$ws.Range("J45").Value = 4484.8335
$ws.Range("K45").Value = 8175.647
$ws.Range("L45").Value = 4484.8335
$ws.Range("M45").Value = -7798.647
$ws.Range("N45").Value = -5238.8335
$ws.Range("H61").Value = 3637.4363
$ws.Range("I61").Value = 3501.1924
$ws.Range("J61").Value = 5999
$ws.Range("K61").Value = 3501.1924
$ws.Range("L61").Value = 5999
$ws.Range("M61").Value = -3289.1924
$ws.Range("N61").Value = -6423
$ws.Range("H63").Value = 2472.55
$ws.Range("I63").Value = 1509.5714
$ws.Range("J63").Value = 2991.077
$ws.Range("K63").Value = 1509.5714
$ws.Range("L63").Value = 2991.077
$ws.Range("M63").Value = -823.5714
$ws.Range("N63").Value = -4363.077
$ws.Range("H66").Value = 2472.55
$ws.Range("I66").Value = 1509.5714
$ws.Range("J66").Value = 2991.077
$ws.Range("K66").Value = 7547.857
$ws.Range("L66").Value = 14955.385
$ws.Range("M66").Value = -4115.857
$ws.Range("N66").Value = -21819.385
$ws.Range("H74").Value = 1742.5098
$ws.Range("I74").Value = 1476.4375
$ws.Range("K74").Value = 1476.4375
$ws.Range("M74").Value = -602.4375
$ws.Range("H77").Value = 1742.5098
$ws.Range("I77").Value = 1476.4375
$ws.Range("K77").Value = 7382.1875
$ws.Range("M77").Value = -3014.1875
$ws.Range("H97").Value = 40041420
$ws.Range("I97").Value = 66668560
$ws.Range("K97").Value = 66668560
$ws.Range("M97").Value = -66668064
$ws.Range("H102").Value = 4831.387
$ws.Range("I102").Value = 2940.1035
$ws.Range("J102").Value = 32255
$ws.Range("K102").Value = 2940.1035
$ws.Range("L102").Value = 32255
$ws.Range("M102").Value = -1318.1035
$ws.Range("N102").Value = -35499
$ws.Range("H122").Value = 8601.578
$ws.Range("I122").Value = 652.69446
$ws.Range("J122").Value = 40397.11
$ws.Range("K122").Value = 1958.08338
$ws.Range("L122").Value = 121191.33
$ws.Range("M122").Value = 491.91662
$ws.Range("N122").Value = -126091.33
$ws.Range("H132").Value = 1680.7949
$ws.Range("I132").Value = 1513.9117
$ws.Range("J132").Value = 2815.6
$ws.Range("K132").Value = 4541.7351
$ws.Range("L132").Value = 8446.799999999999
$ws.Range("M132").Value = -2011.7351
$ws.Range("N132").Value = -13506.8
$ws.Range("H136").Value = 3637.4363
$ws.Range("I136").Value = 3501.1924
$ws.Range("J136").Value = 5999
$ws.Range("K136").Value = 10503.5772
$ws.Range("L136").Value = 17997
$ws.Range("M136").Value = -7953.5772
$ws.Range("N136").Value = -23097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3575.1365
$ws.Range("I20").Value = 3242.7036
$ws.Range("J20").Value = 4103.1177
$ws.Range("K20").Value = 3242.7036
$ws.Range("L20").Value = 4103.1177
$ws.Range("M20").Value = -2995.7036
$ws.Range("N20").Value = -4597.1177
$ws.Range("H94").Value = 382.87878
$ws.Range("I94").Value = 312.12
$ws.Range("J94").Value = 604
$ws.Range("K94").Value = 312.12
$ws.Range("L94").Value = 604
$ws.Range("M94").Value = 138.88
$ws.Range("N94").Value = -1506
$ws.Range("H99").Value = 1929.9546
$ws.Range("I99").Value = 1719.3334
$ws.Range("J99").Value = 2877.75
$ws.Range("K99").Value = 1719.3334
$ws.Range("L99").Value = 2877.75
$ws.Range("M99").Value = -221.3334
$ws.Range("N99").Value = -5873.75
$ws.Range("H105").Value = 2725.5
$ws.Range("I105").Value = 1300.1154
$ws.Range("J105").Value = 21255.5
$ws.Range("K105").Value = 1300.1154
$ws.Range("L105").Value = 21255.5
$ws.Range("M105").Value = 446.8846000000001
$ws.Range("N105").Value = -24749.5
$ws.Range("H134").Value = 3108.0715
$ws.Range("I134").Value = 3070.5386
$ws.Range("J134").Value = 3596
$ws.Range("K134").Value = 9211.6158
$ws.Range("L134").Value = 10788
$ws.Range("M134").Value = -6676.6158
$ws.Range("N134").Value = -15858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 994.5
$ws.Range("I22").Value = 990
$ws.Range("J22").Value = 999
$ws.Range("K22").Value = 990
$ws.Range("L22").Value = 999
$ws.Range("M22").Value = -640
$ws.Range("N22").Value = -1699
$ws.Range("H58").Value = 2010.5135
$ws.Range("I58").Value = 1955.7
$ws.Range("J58").Value = 2245.4285
$ws.Range("K58").Value = 1955.7
$ws.Range("L58").Value = 2245.4285
$ws.Range("M58").Value = -1752.7
$ws.Range("N58").Value = -2651.4285
$ws.Range("H136").Value = 2010.5135
$ws.Range("I136").Value = 1955.7
$ws.Range("J136").Value = 2245.4285
$ws.Range("K136").Value = 5867.1
$ws.Range("L136").Value = 6736.2855
$ws.Range("M136").Value = -3317.1
$ws.Range("N136").Value = -11836.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 36925844
$ws.Range("I80").Value = 68573710
$ws.Range("J80").Value = 3332.5
$ws.Range("K80").Value = 68573710
$ws.Range("L80").Value = 3332.5
$ws.Range("M80").Value = -68572712
$ws.Range("N80").Value = -5328.5
$ws.Range("H83").Value = 36925844
$ws.Range("I83").Value = 68573710
$ws.Range("J83").Value = 3332.5
$ws.Range("K83").Value = 342868550
$ws.Range("L83").Value = 16662.5
$ws.Range("M83").Value = -342863558
$ws.Range("N83").Value = -26646.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2522.7693
$ws.Range("I7").Value = 2131
$ws.Range("J7").Value = 3149.6
$ws.Range("K7").Value = 2131
$ws.Range("L7").Value = 3149.6
$ws.Range("M7").Value = -2019
$ws.Range("N7").Value = -3373.6
$ws.Range("H16").Value = 8929637
$ws.Range("I16").Value = 12500860
$ws.Range("J16").Value = 1578
$ws.Range("K16").Value = 12500860
$ws.Range("L16").Value = 1578
$ws.Range("M16").Value = -12500690
$ws.Range("N16").Value = -1918
$ws.Range("H22").Value = 855.9091
$ws.Range("I22").Value = 846.6667
$ws.Range("J22").Value = 897.5
$ws.Range("K22").Value = 846.6667
$ws.Range("L22").Value = 897.5
$ws.Range("M22").Value = -551.6667
$ws.Range("N22").Value = -1487.5
$ws.Range("H27").Value = 855.9091
$ws.Range("I27").Value = 846.6667
$ws.Range("J27").Value = 897.5
$ws.Range("K27").Value = 846.6667
$ws.Range("L27").Value = 897.5
$ws.Range("M27").Value = -739.6667
$ws.Range("N27").Value = -1111.5
$ws.Range("H61").Value = 2471.7
$ws.Range("I61").Value = 2675.7144
$ws.Range("J61").Value = 1995.6666
$ws.Range("K61").Value = 2675.7144
$ws.Range("L61").Value = 1995.6666
$ws.Range("M61").Value = -2473.7144
$ws.Range("N61").Value = -2399.6666
$ws.Range("H82").Value = 1993.5
$ws.Range("I82").Value = 2062.8708
$ws.Range("J82").Value = 1071.8572
$ws.Range("K82").Value = 2062.8708
$ws.Range("L82").Value = 1071.8572
$ws.Range("M82").Value = -1701.8708
$ws.Range("N82").Value = -1793.8572
$ws.Range("H85").Value = 1993.5
$ws.Range("I85").Value = 2062.8708
$ws.Range("J85").Value = 1071.8572
$ws.Range("K85").Value = 2062.8708
$ws.Range("L85").Value = 1071.8572
$ws.Range("M85").Value = -814.8708000000001
$ws.Range("N85").Value = -3567.8572
$ws.Range("H100").Value = 71431280
$ws.Range("I100").Value = 250002080
$ws.Range("J100").Value = 2958.9
$ws.Range("K100").Value = 250002080
$ws.Range("L100").Value = 2958.9
$ws.Range("M100").Value = -250001539
$ws.Range("N100").Value = -4040.9
$ws.Range("H113").Value = 2471.7
$ws.Range("I113").Value = 2675.7144
$ws.Range("J113").Value = 1995.6666
$ws.Range("K113").Value = 2675.7144
$ws.Range("L113").Value = 1995.6666
$ws.Range("M113").Value = -505.7143999999998
$ws.Range("N113").Value = -6335.6666
$ws.Range("H126").Value = 2522.7693
$ws.Range("I126").Value = 2131
$ws.Range("J126").Value = 3149.6
$ws.Range("K126").Value = 6393
$ws.Range("L126").Value = 9448.799999999999
$ws.Range("M126").Value = -3923
$ws.Range("N126").Value = -14388.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 602.3889
$ws.Range("I100").Value = 529.5
$ws.Range("J100").Value = 748.1667
$ws.Range("K100").Value = 1059
$ws.Range("L100").Value = 1496.3334
$ws.Range("M100").Value = -518
$ws.Range("N100").Value = -2578.3334
$ws.Range("H107").Value = 1399.1578
$ws.Range("I107").Value = 1480.1111
$ws.Range("J107").Value = 1326.3
$ws.Range("K107").Value = 4440.3333
$ws.Range("L107").Value = 3978.9
$ws.Range("M107").Value = -2520.3333
$ws.Range("N107").Value = -7818.9
$ws.Range("H113").Value = 1184.55
$ws.Range("I113").Value = 1105.4706
$ws.Range("J113").Value = 1632.6666
$ws.Range("K113").Value = 3316.4118
$ws.Range("L113").Value = 4897.9998
$ws.Range("M113").Value = -1146.4118
$ws.Range("N113").Value = -9237.9998
$ws.Range("H122").Value = 3468.0364
$ws.Range("I122").Value = 1995.279
$ws.Range("J122").Value = 8745.416999999999
$ws.Range("K122").Value = 5985.837
$ws.Range("L122").Value = 26236.251
$ws.Range("M122").Value = -3535.837
$ws.Range("N122").Value = -31136.251
$ws.Range("H132").Value = 2597.4827
$ws.Range("I132").Value = 1810.1818
$ws.Range("J132").Value = 5071.857
$ws.Range("K132").Value = 5430.5454
$ws.Range("L132").Value = 15215.571
$ws.Range("M132").Value = -2900.5454
$ws.Range("N132").Value = -20275.571
$ws.Range("H136").Value = 2573.1633
$ws.Range("I136").Value = 2188.6047
$ws.Range("J136").Value = 5329.1665
$ws.Range("K136").Value = 6565.8141
$ws.Range("L136").Value = 15987.4995
$ws.Range("M136").Value = -7329.1665
$ws.Range("N136").Value = -21087.4995
